# Add new daily-routine rows (21-29) to the log sheet, matching the
# author's upload of additional progress entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column values (row -> serial date number) and whether the row
# uses the alternate "d-mmm" number format (only row 26 in the diff).
$rows = @(
    @{ Row = 21; Date = 43853; Text = "create the question randomize on the level's"; AltFormat = $false },
    @{ Row = 22; Date = 43854; Text = "design the home page of the take the test";    AltFormat = $false },
    @{ Row = 23; Date = 43857; Text = "design the admin page";                        AltFormat = $false },
    @{ Row = 24; Date = 43858; Text = "insert the database for sample data's";        AltFormat = $false },
    @{ Row = 25; Date = 43859; Text = "connection into the java";                     AltFormat = $false },
    @{ Row = 26; Date = 43860; Text = "insert the database for sample data's in jsp"; AltFormat = $true  },
    @{ Row = 27; Date = 43861; Text = "question to insert the html page to database"; AltFormat = $false },
    @{ Row = 28; Date = 43865; Text = "document correction";                         AltFormat = $false },
    @{ Row = 29; Date = 43866; Text = "split the question";                          AltFormat = $false }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $dateCell = $ws.Cells.Item($rowNum, 1)
    $dateCell.Value = $r.Date
    if ($r.AltFormat) {
        $dateCell.NumberFormat = "d-mmm"
    } else {
        $dateCell.NumberFormat = "m/d/yy"
    }

    $ws.Cells.Item($rowNum, 2).Value = $r.Text
}

# Update the view to reflect where the user ended up after the edit.
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B29").Select() | Out-Null
